$d = $word.ActiveDocument

# --- 1. Locate the anchor sentence ("En måde at begrænse kode på og gøre
#        det mere overskueligt. ") that ends the "Funktioner" paragraph,
#        and split the paragraph right after it (twice), so we end up
#        with a brand new, empty paragraph sitting between the anchor
#        sentence and the old bookmark/page-break paragraph. ---
$find = $d.Content
$found = $find.Find.Execute("En måde at begrænse kode på og gøre det mere overskueligt. ", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

$sentStart = $find.Start
$sentEnd = $find.End

$anchor = $d.Range($sentStart, $sentEnd)
$anchor.InsertParagraphAfter()

$anchor2 = $d.Range($sentStart, $sentEnd)
$anchor2.InsertParagraphAfter()

# --- 2. Fill the freshly opened, empty paragraph with the new content. ---
$pos = $sentEnd + 1

$ins = $d.Range($pos, $pos)
$ins.InsertAfter("Man kan bruge ")
$pos = $pos + "Man kan bruge ".Length

$t = "p"
$ins = $d.Range($pos, $pos)
$ins.InsertAfter($t)
$r2 = $d.Range($pos, $pos + $t.Length)
$r2.Font.Italic = $true
$pos = $pos + $t.Length

$t = "arameter"
$ins = $d.Range($pos, $pos)
$ins.InsertAfter($t)
$r2 = $d.Range($pos, $pos + $t.Length)
$r2.Font.Italic = $true
$pos = $pos + $t.Length

$ins = $d.Range($pos, $pos)
$ins.InsertAfter(" i funktioner fx ")
$pos = $pos + " i funktioner fx ".Length

$ins = $d.Range($pos, $pos)
$ins.InsertAfter("function")
$pos = $pos + "function".Length

$ins = $d.Range($pos, $pos)
$ins.InsertAfter(" ")
$pos = $pos + " ".Length

$ins = $d.Range($pos, $pos)
$ins.InsertAfter("findSum")
$pos = $pos + "findSum".Length

$ins = $d.Range($pos, $pos)
$ins.InsertAfter("(")
$pos = $pos + "(".Length

$ins = $d.Range($pos, $pos)
$ins.InsertAfter("a, b)")
$pos = $pos + "a, b)".Length

$ins = $d.Range($pos, $pos)
$ins.InsertAfter(" hvor a og b er ")
$pos = $pos + " hvor a og b er ".Length

$ins = $d.Range($pos, $pos)
$ins.InsertAfter("parametrene")
$pos = $pos + "parametrene".Length

$ins = $d.Range($pos, $pos)
$ins.InsertAfter(",")
$pos = $pos + ",".Length

$ins = $d.Range($pos, $pos)
$ins.InsertAfter(" som man så kan kalde ved at skrive ")
$pos = $pos + " som man så kan kalde ved at skrive ".Length

$t = "argumenterne"
$ins = $d.Range($pos, $pos)
$ins.InsertAfter($t)
$r2 = $d.Range($pos, $pos + $t.Length)
$r2.Font.Italic = $true
$pos = $pos + $t.Length

$ins = $d.Range($pos, $pos)
$ins.InsertAfter(" i funktionen ")
$pos = $pos + " i funktionen ".Length

$ins = $d.Range($pos, $pos)
$ins.InsertAfter("fx ")
$pos = $pos + "fx ".Length

$ins = $d.Range($pos, $pos)
$ins.InsertAfter("findSum")
$pos = $pos + "findSum".Length

$ins = $d.Range($pos, $pos)
$ins.InsertAfter("(8, 90); ")
$pos = $pos + "(8, 90); ".Length
